$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New HIGH / LOW / CLOSE / LTP / VOL / 9:25 CLOSE values (columns B-G) for rows 2-17
# fixing a date/month lookup bug in sonia.py.
$data = @{
    2  = @(814,       788.75,   805,      804.85,   98,   790.5)
    3  = @(842,       825.5,    828.25,   828.35,   26,   837.75)
    4  = @(44746.15,  44061,    44630,    44601.5,  34,   44263.1)
    5  = @(332.6,     322.15,   330.8,    330.85,   233,  322.85)
    6  = @(512.45,    502.6,    505.95,   505.75,   73,   507.3)
    7  = @(481.7,     472.05,   476.1,    476.8,    181,  472.7)
    8  = @(976.75,    958.5,    971.35,   972.05,   174,  962.2)
    9  = @(719,       694.85,   704.9,    705.8,    113,  696.5)
    10 = @(19528.95,  19380,    19519.5,  19512.1,  58,   19390)
    11 = @(2440.45,   2419.05,  2427,     2428,     85,   2424.9)
    12 = @(573.25,    564.9,    571.65,   571.4,    227,  567.15)
    13 = @(851.5,     843,      846.05,   849.8,    15,   843.2)
    14 = @(619.2,     607,      612.9,    613.95,   223,  607.3)
    15 = @(129.5,     126.3,    127.8,    128,      1008, 126.65)
    16 = @(3404,      3385.25,  3396.1,   3400.1,   12,   3388.8)
    17 = @(3138,      3112.05,  3129.05,  3128.1,   6,    3132.8)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
